$wb = $excel.ActiveWorkbook

# The handback report was generated for the c65c4556-... file (row 8) in
# both locale sheets (zh-cn / de-de): its "Latest Target File", "Latest
# Handback File", "Latest Handback DateTime" and "Error Detail" columns
# get populated, and the I/P columns widen to fit the new content.

$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/76e72f6a7e76bd01f43d7f2a47b23c47af09a184/e2e/c65c4556-3114-4556-8f29-6575ac282f18.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/14dbe1c177dcb7b78e2407824be7762971bba81d/e2e/c65c4556-3114-4556-8f29-6575ac282f18.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/76e72f6a7e76bd01f43d7f2a47b23c47af09a184/e2e/c65c4556-3114-4556-8f29-6575ac282f18.md."

function Update-LocaleSheet {
    param($SheetName, $HandbackFile, $HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # I8 - Latest Target File: becomes a hyperlink to the handed-back md file
    $ws.Range("I8").Value = "c65c4556-3114-4556-8f29-6575ac282f18.md"
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetMdUrl, $null, $null, "c65c4556-3114-4556-8f29-6575ac282f18.md")
    $ws.Range("I8").Font.Color = 15570276
    $ws.Range("I8").Font.Underline = 2

    # J8 - Latest Handback File
    $ws.Range("J8").Value = $HandbackFile

    # K8 - Latest Handback DateTime
    $ws.Range("K8").Value = $HandbackDateTime

    # P8 - Error Detail
    $ws.Range("P8").Value = $errorDetail

    # Widen columns I and P so the new content fits (both become 40 chars).
    $ws.Columns.Item(9).ColumnWidth = 39.14
    $ws.Columns.Item(16).ColumnWidth = 39.14
}

Update-LocaleSheet "zh-cn" "c65c4556-3114-4556-8f29-6575ac282f18.68d914556a21a627e27bf9289829a54d447f9ad1.zh-cn.xlf" "2016-09-05 00:49:14"
Update-LocaleSheet "de-de" "c65c4556-3114-4556-8f29-6575ac282f18.68d914556a21a627e27bf9289829a54d447f9ad1.de-de.xlf" "2016-09-05 00:49:21"
